# Auto-generated Excel COM-interop script applying the scheduled-runner data refresh
# described by the commit diff (per-sheet "Leve profit" recompute).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 946
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 1100
$ws.Range("I16").Value = 1100
$ws.Range("K16").Value = 1100
$ws.Range("M16").Value = -813
$ws.Range("H45").Value = 3843.4
$ws.Range("J45").Value = 6406.6665
$ws.Range("L45").Value = 6406.6665
$ws.Range("N45").Value = -7160.6665
$ws.Range("H61").Value = 15551.75
$ws.Range("I61").Value = 3150
$ws.Range("J61").Value = 19685.666
$ws.Range("K61").Value = 3150
$ws.Range("L61").Value = 19685.666
$ws.Range("M61").Value = -2938
$ws.Range("N61").Value = -20109.666
$ws.Range("H63").Value = 2073
$ws.Range("I63").Value = 2073
$ws.Range("K63").Value = 2073
$ws.Range("M63").Value = -1387
$ws.Range("H66").Value = 2073
$ws.Range("I66").Value = 2073
$ws.Range("K66").Value = 10365
$ws.Range("M66").Value = -6933
$ws.Range("H74").Value = 471138.16
$ws.Range("I74").Value = 1200712.2
$ws.Range("J74").Value = 15154.375
$ws.Range("K74").Value = 1200712.2
$ws.Range("L74").Value = 15154.375
$ws.Range("M74").Value = -1199838.2
$ws.Range("N74").Value = -16902.375
$ws.Range("H77").Value = 471138.16
$ws.Range("I77").Value = 1200712.2
$ws.Range("J77").Value = 15154.375
$ws.Range("K77").Value = 6003561
$ws.Range("L77").Value = 75771.875
$ws.Range("M77").Value = -5999193
$ws.Range("N77").Value = -84507.875
$ws.Range("H88").Value = 4601.1333
$ws.Range("I88").Value = 1202.1
$ws.Range("J88").Value = 11399.2
$ws.Range("K88").Value = 1202.1
$ws.Range("L88").Value = 11399.2
$ws.Range("M88").Value = -796.0999999999999
$ws.Range("N88").Value = -12211.2
$ws.Range("H91").Value = 4601.1333
$ws.Range("I91").Value = 1202.1
$ws.Range("J91").Value = 11399.2
$ws.Range("K91").Value = 1202.1
$ws.Range("L91").Value = 11399.2
$ws.Range("M91").Value = 201.9000000000001
$ws.Range("N91").Value = -14207.2
$ws.Range("H110").Value = 9574.5
$ws.Range("I110").Value = 9574.5
$ws.Range("K110").Value = 9574.5
$ws.Range("M110").Value = -7529.5
$ws.Range("H132").Value = 2017.6
$ws.Range("I132").Value = 1852.8889
$ws.Range("J132").Value = 3500
$ws.Range("K132").Value = 5558.6667
$ws.Range("L132").Value = 10500
$ws.Range("M132").Value = -3028.6667
$ws.Range("N132").Value = -15560
$ws.Range("H136").Value = 15551.75
$ws.Range("I136").Value = 3150
$ws.Range("J136").Value = 19685.666
$ws.Range("K136").Value = 9450
$ws.Range("L136").Value = 59056.99800000001
$ws.Range("M136").Value = -6900
$ws.Range("N136").Value = -64156.99800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 21812.875
$ws.Range("H85").Value = 21812.875
$ws.Range("H86").Value = 1479.4166
$ws.Range("I86").Value = 1044.5454
$ws.Range("K86").Value = 1044.5454
$ws.Range("M86").Value = 78.45460000000003
$ws.Range("H89").Value = 1479.4166
$ws.Range("I89").Value = 1044.5454
$ws.Range("K89").Value = 5222.727
$ws.Range("M89").Value = 393.2730000000001
$ws.Range("H94").Value = 1305
$ws.Range("I94").Value = 599.05884
$ws.Range("K94").Value = 599.05884
$ws.Range("M94").Value = -148.05884
$ws.Range("H99").Value = 1526.4445
$ws.Range("I99").Value = 1217.25
$ws.Range("J99").Value = 4000
$ws.Range("K99").Value = 1217.25
$ws.Range("L99").Value = 4000
$ws.Range("M99").Value = 280.75
$ws.Range("N99").Value = -6996
$ws.Range("H105").Value = 2744.6
$ws.Range("J105").Value = 4117.6665
$ws.Range("L105").Value = 4117.6665
$ws.Range("N105").Value = -7611.6665
$ws.Range("H107").Value = 2953.7368
$ws.Range("I107").Value = 2919.3333
$ws.Range("J107").Value = 3012.7144
$ws.Range("K107").Value = 2919.3333
$ws.Range("L107").Value = 3012.7144
$ws.Range("M107").Value = -999.3332999999998
$ws.Range("N107").Value = -6852.7144
$ws.Range("H134").Value = 2717.5938
$ws.Range("I134").Value = 2618.138
$ws.Range("J134").Value = 3679
$ws.Range("K134").Value = 7854.414
$ws.Range("L134").Value = 11037
$ws.Range("M134").Value = -5319.414
$ws.Range("N134").Value = -16107

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 39999
$ws.Range("J18").Value = 39999
$ws.Range("L18").Value = 39999
$ws.Range("N18").Value = -40459
$ws.Range("H31").Value = 5004099.5
$ws.Range("I31").Value = 25000750
$ws.Range("J31").Value = 4937.0625
$ws.Range("K31").Value = 25000750
$ws.Range("L31").Value = 4937.0625
$ws.Range("M31").Value = -25000455
$ws.Range("N31").Value = -5527.0625
$ws.Range("H34").Value = 5004099.5
$ws.Range("I34").Value = 25000750
$ws.Range("J34").Value = 4937.0625
$ws.Range("K34").Value = 25000750
$ws.Range("L34").Value = 4937.0625
$ws.Range("M34").Value = -25000548
$ws.Range("N34").Value = -5341.0625
$ws.Range("H58").Value = 1022.68964
$ws.Range("I58").Value = 860.44446
$ws.Range("J58").Value = 3213
$ws.Range("K58").Value = 860.44446
$ws.Range("L58").Value = 3213
$ws.Range("M58").Value = -657.44446
$ws.Range("N58").Value = -3619
$ws.Range("H136").Value = 1022.68964
$ws.Range("I136").Value = 860.44446
$ws.Range("J136").Value = 3213
$ws.Range("K136").Value = 2581.33338
$ws.Range("L136").Value = 9639
$ws.Range("M136").Value = -31.33338000000003
$ws.Range("N136").Value = -14739

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 1212.8
$ws.Range("I3").Value = 1212.8
$ws.Range("K3").Value = 3638.4
$ws.Range("M3").Value = -3526.4
$ws.Range("H5").Value = 755.9091
$ws.Range("I5").Value = 741.5
$ws.Range("J5").Value = 900
$ws.Range("K5").Value = 2224.5
$ws.Range("L5").Value = 2700
$ws.Range("M5").Value = -2112.5
$ws.Range("N5").Value = -2924
$ws.Range("H46").Value = 910.8889
$ws.Range("I46").Value = 900
$ws.Range("J46").Value = 932.6667
$ws.Range("K46").Value = 2700
$ws.Range("L46").Value = 2798.0001
$ws.Range("M46").Value = -2609
$ws.Range("N46").Value = -2980.0001
$ws.Range("H59").Value = 5958.759
$ws.Range("J59").Value = 4307.654
$ws.Range("L59").Value = 12922.962
$ws.Range("N59").Value = -14002.962
$ws.Range("H68").Value = 4595.696
$ws.Range("J68").Value = 4768.136
$ws.Range("L68").Value = 14304.408
$ws.Range("N68").Value = -15926.408
$ws.Range("H71").Value = 4595.696
$ws.Range("J71").Value = 4768.136
$ws.Range("L71").Value = 42913.224
$ws.Range("N71").Value = -51025.224
$ws.Range("H135").Value = 755.9091
$ws.Range("I135").Value = 741.5
$ws.Range("J135").Value = 900
$ws.Range("K135").Value = 6673.5
$ws.Range("L135").Value = 8100
$ws.Range("M135").Value = -4138.5
$ws.Range("N135").Value = -13170

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 21705.572
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 21705.572
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 21705.572
$ws.Range("M24").ClearContents()
$ws.Range("N24").Value = -22051.572
$ws.Range("H122").Value = 5342.75
$ws.Range("I122").Value = 5229.615
$ws.Range("J122").Value = 5833
$ws.Range("K122").Value = 15688.845
$ws.Range("L122").Value = 17499
$ws.Range("M122").Value = -13238.845
$ws.Range("N122").Value = -22399
$ws.Range("H132").Value = 2154.8147
$ws.Range("I132").Value = 1984.762
$ws.Range("J132").Value = 2750
$ws.Range("K132").Value = 5954.286
$ws.Range("L132").Value = 8250
$ws.Range("M132").Value = -3424.286
$ws.Range("N132").Value = -13310

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2942.889
$ws.Range("I61").Value = 2087.7778
$ws.Range("J61").Value = 3798
$ws.Range("K61").Value = 2087.7778
$ws.Range("L61").Value = 3798
$ws.Range("M61").Value = -1885.7778
$ws.Range("N61").Value = -4202
$ws.Range("H63").Value = 75000
$ws.Range("J63").Value = 75000
$ws.Range("L63").Value = 75000
$ws.Range("N63").Value = -76498
$ws.Range("H66").Value = 75000
$ws.Range("J66").Value = 75000
$ws.Range("L66").Value = 225000
$ws.Range("N66").Value = -232488
$ws.Range("H68").Value = 3986.889
$ws.Range("I68").Value = 3926
$ws.Range("K68").Value = 3926
$ws.Range("M68").Value = -3177
$ws.Range("H71").Value = 3986.889
$ws.Range("I71").Value = 3926
$ws.Range("K71").Value = 19630
$ws.Range("M71").Value = -15886
$ws.Range("H75").Value = 25000
$ws.Range("J75").Value = 25000
$ws.Range("L75").Value = 25000
$ws.Range("N75").Value = -26872
$ws.Range("H78").Value = 25000
$ws.Range("J78").Value = 25000
$ws.Range("L78").Value = 75000
$ws.Range("N78").Value = -84360
$ws.Range("H113").Value = 2942.889
$ws.Range("I113").Value = 2087.7778
$ws.Range("J113").Value = 3798
$ws.Range("K113").Value = 2087.7778
$ws.Range("L113").Value = 3798
$ws.Range("M113").Value = 82.22220000000016
$ws.Range("N113").Value = -8138
$ws.Range("H122").Value = 3351.1765
$ws.Range("J122").Value = 3097.1
$ws.Range("L122").Value = 9291.299999999999
$ws.Range("N122").Value = -14191.3
$ws.Range("H136").Value = 3053.7666
$ws.Range("I136").Value = 2680.2
$ws.Range("K136").Value = 8040.599999999999
$ws.Range("M136").Value = -5490.599999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 2115.3845
$ws.Range("I2").Value = 2226.6667
$ws.Range("J2").Value = 780
$ws.Range("K2").Value = 2226.6667
$ws.Range("L2").Value = 780
$ws.Range("M2").Value = -2114.6667
$ws.Range("N2").Value = -1004
$ws.Range("H122").Value = 32984.14
$ws.Range("I122").Value = 37893.918
$ws.Range("J122").Value = 2707.1667
$ws.Range("K122").Value = 113681.754
$ws.Range("L122").Value = 8121.500100000001
$ws.Range("M122").Value = -111231.754
$ws.Range("N122").Value = -13021.5001
$ws.Range("H132").Value = 18837.121
$ws.Range("I132").Value = 21872.795
$ws.Range("K132").Value = 65618.38499999999
$ws.Range("M132").Value = -63088.38499999999
